$d = $word.ActiveDocument

# Step 1: Append the new sentence run after "...to generate the data. "
$r1 = $d.Content
$found1 = $r1.Find.Execute(", the larger it is, the more files we Java needs to open to generate the data. ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$insertPos = $r1.End
$appendRange = $d.Range($insertPos, $insertPos)
$appendRange.InsertAfter("Of course, since we need to generate the user tree, it’s instantly a linear operation. Performing queries for this would end up being in quasilinear. On an unrelated note, since we changed from C to Java between stage_2_c and stage_3, the time taken will for sure be different due to the limitations and performance of Java. ")

# Step 2: Move the _GoBack bookmark to split "nodes" into "n" | "odes"
#         (Bookmarks.Add with the existing name relocates it, removing the old one.)
$r2 = $d.Content
$found2 = $r2.Find.Execute(" Since we have more n", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$splitPos = $r2.End
$splitRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $splitRange)
